$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("27072019")

# Insert a new blank row at the top, shifting the header + data down by one row.
$ws.Rows.Item(1).Insert()

# Update the (now shifted) header row (row 2) text to match the "Cell Cycle" sheet's
# header wording/shared strings.
$ws.Range("B2").Value = " %G1"
$ws.Range("C2").Value = " %S"
$ws.Range("D2").Value = " %G2"
$ws.Range("E2").Value = " G1 Mean"
$ws.Range("F2").Value = "G2 Mean"
$ws.Range("G2").Value = " G1 CV"
$ws.Range("H2").Value = " G2 CV"
$ws.Range("I2").Value = " % < G1"
$ws.Range("J2").Value = " % > G2"

# Re-apply the AutoFilter over the new header row (A2:J2).
$ws.AutoFilterMode = $false
$ws.Range("A2:J2").AutoFilter()

# Keep the workbook-level _FilterDatabase defined name in sync with the new range.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -eq "27072019!_FilterDatabase") {
        $n.RefersTo = "='27072019'!`$A`$2:`$J`$2"
    }
}

# Restore the active selection.
$ws.Range("D11").Select()
